$d = $word.ActiveDocument

# --- Paragraph 1: the **ID__...__ID** marker paragraph ---
$p1 = $d.Paragraphs(1)

# Replace the old topic id placeholder text with the new one.
$null = $p1.Range.Find.Execute(
    "**ID__AFFARS_pgi_5306_topic_11__ID**", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5306_502__ID**", 2)

# The paragraph used to have a second run containing just a trailing
# space; drop that trailing space so the paragraph text ends right
# after "__ID**" (the now-merged run no longer needs xml:space="preserve").
$null = $p1.Range.Find.Execute(
    "**ID__AFFARS_SMC_PGI_5306_502__ID** ", $true, $false, $false, $false,
    $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5306_502__ID**", 2)

# Re-fetch the paragraph (defensive -- content changed) and give it the
# same paragraph-border / indent treatment used elsewhere in this
# document: a 5-twip spaced border box and a 225-twip (11.25pt) left
# indent instead of the previous 120-twip (6pt) indent.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
